$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.636.33'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '2.308.10'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '268.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.620'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.60'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0933'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.31%  '
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').Value = '2.654.28'
$ws.Range('E14').Value = '  +1.76%  '
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.863'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.03%  '
$ws.Range('D17').Value = '2.312.18'
$ws.Range('E17').Value = '  +2.38%  '
$ws.Range('D18').Value = '43.629.80'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('E19').Value = '  +1.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.45%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.50'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.94%  '
$ws.Range('E29').Value = '  -4.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.62'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '171.64'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0896'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.126'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0359'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.48'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.108'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.233'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +14.09%  '
$ws.Range('E41').Value = '  +3.12%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +17.13%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.94%  '
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '100.24'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.84%  '
$ws.Range('E49').Value = '  -1.45%  '
$ws.Range('D50').Value = '2.532.46'
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.426'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.61%  '
